$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the id values in column A for rows 2-4
$ws.Range("A2").Value = 611
$ws.Range("A3").Value = 607
$ws.Range("A4").Value = 608

# Move the active selection from A5 to A2
$ws.Range("A2").Select()
